$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# Update status text for the e1dd92d3 row (row 7) on all sheets that show it.
$overview.Range("B7").Value = "Handback transform failed"
$overview.Range("C7").Value = "Handback transform failed"
$zhcn.Range("C7").Value = "Handback transform failed"
$dede.Range("C7").Value = "Handback transform failed"

# Add "Error Detail" (column L) values for row 7 on the locale sheets.
$zhcn.Range("L7").Value = "Handback file name: rwgn4yuq.vut is different with handoff file name: e1dd92d3-efea-4269-bf79-2b79ca59b586.db5ebfea663117e9c389f347f8b496b5ea63f10d.zh-cn."
$dede.Range("L7").Value = "Handback file name: rwgn4yuq.vut is different with handoff file name: e1dd92d3-efea-4269-bf79-2b79ca59b586.db5ebfea663117e9c389f347f8b496b5ea63f10d.de-de."
